$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text assignments (values that Excel will not misinterpret as numbers)
$directValues = @{
    'D2' = '57.129.49'
    'E2' = '  -2.27%  '
    'D3' = '3.067.52'
    'E3' = '  -2.39%  '
    'E4' = '  +0.00%  '
    'E5' = '  -2.45%  '
    'E6' = '  -5.45%  '
    'E7' = '  +0.03%  '
    'D8' = '3.066.06'
    'E8' = '  -2.42%  '
    'E9' = '  +5.67%  '
    'E10' = '  +0.78%  '
    'E11' = '  -3.49%  '
    'E12' = '  +1.23%  '
    'E13' = '  +1.68%  '
    'D14' = '3.596.06'
    'E14' = '  -2.41%  '
    'E15' = '  -2.77%  '
    'E16' = '  -4.49%  '
    'D17' = '57.158.96'
    'E17' = '  -2.28%  '
    'D18' = '3.066.76'
    'E18' = '  -2.77%  '
    'E19' = '  -4.61%  '
    'E20' = '  -3.90%  '
    'E21' = '  -2.96%  '
    'E22' = '  +0.63%  '
    'E23' = '  -0.21%  '
    'E24' = '  +1.52%  '
    'E25' = '  -3.33%  '
    'E26' = '  +0.12%  '
    'E27' = '  -2.83%  '
    'D28' = '0.0₃0836'
    'E28' = '  -10.43%  '
    'E29' = '  -0.03%  '
    'E30' = '  -5.57%  '
    'E31' = '  -3.22%  '
    'E32' = '  -1.26%  '
    'E33' = '  -10.43%  '
    'E34' = '  +0.09%  '
    'E35' = '  -0.80%  '
    'E36' = '  -7.35%  '
    'E37' = '  -5.00%  '
    'E38' = '  -4.29%  '
    'E39' = '  -3.45%  '
    'E40' = '  -3.11%  '
    'E41' = '  -5.83%  '
    'E42' = '  -0.54%  '
    'E43' = '  -3.10%  '
    'D44' = '2.403.91'
    'E44' = '  +4.75%  '
    'E45' = '  -0.36%  '
    'D47' = '3.108.55'
    'E47' = '  -2.38%  '
    'E48' = '  -2.10%  '
    'E49' = '  -2.68%  '
    'E50' = '  -8.57%  '
    'E51' = '  -7.18%  '
}
foreach ($cell in $directValues.Keys) {
    $ws.Range($cell).Value = $directValues[$cell]
}

# Values that look numeric ("520.14", "1.00", etc.) must be forced to stay plain text,
# matching the original inline-string cells, without picking up a new NumberFormat/
# style. Build each one as a text formula in a scratch cell, then paste-special just
# the value back onto the target cell so no style/formula residue is left behind.
$numericTextValues = @{
    'D5' = '520.14'
    'D6' = '135.36'
    'D7' = '1.00'
    'D9' = '0.472'
    'D10' = '7.23'
    'D12' = '0.399'
    'D15' = '24.99'
    'D22' = '346.74'
    'D23' = '0.999'
    'D24' = '68.84'
    'D25' = '0.496'
    'D31' = '1.83'
    'D32' = '20.89'
    'D33' = '5.77'
    'D34' = '158.23'
    'D36' = '1.11'
    'D37' = '5.95'
    'D38' = '25.17'
    'D39' = '1.21'
    'D40' = '0.0653'
    'D42' = '4.00'
    'D49' = '5.94'
    'D50' = '0.928'
    'D51' = '19.23'
}
$scratch = $ws.Range("ZZ1")
foreach ($cell in $numericTextValues.Keys) {
    $scratch.Formula = '="' + $numericTextValues[$cell] + '"'
    $scratch.Copy()
    $ws.Range($cell).PasteSpecial(-4163)
}
$scratch.Value = ""
$excel.CutCopyMode = 0

